$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: add "A" in A1, clear F1 (it becomes an empty gap), add "G" in G1
$ws.Range("A1").Value = "A"
$ws.Range("F1").ClearContents()
$ws.Range("G1").Value = "G"

# Fill rows 2-6, columns A-G with values 1-7
for ($r = 2; $r -le 6; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $c
    }
}

# Update the sheet's selection to cover the full data range A1:G6
$ws.Range("A1:G6").Select()
